# Replace the three-digit x one-digit multiplication problems/answers
# in the table with the new set of values.

$d = $word.ActiveDocument

$replacements = @(
    @("698×8=5584", "762×4=3048"),
    @("925×4=3700", "638×5=3190"),
    @("299×5=1495", "598×2=1196"),
    @("912×9=8208", "447×7=3129"),
    @("781×7=5467", "974×8=7792"),
    @("221×6=1326", "569×7=3983"),
    @("673×6=4038", "712×5=3560"),
    @("485×9=4365", "681×4=2724"),
    @("546×2=1092", "652×7=4564"),
    @("519×5=2595", "170×5=850"),
    @("756×7=5292", "563×2=1126"),
    @("772×2=1544", "507×4=2028"),
    @("143×5=715",  "562×2=1124"),
    @("980×5=4900", "650×2=1300"),
    @("406×7=2842", "659×9=5931"),
    @("791×5=3955", "707×3=2121"),
    @("610×6=3660", "518×3=1554"),
    @("879×8=7032", "490×8=3920"),
    @("841×2=1682", "470×3=1410"),
    @("663×9=5967", "300×5=1500"),
    @("194×6=1164", "744×8=5952"),
    @("440×6=2640", "753×3=2259"),
    @("217×8=1736", "168×8=1344"),
    @("144×9=1296", "659×7=4613"),
    @("807×2=1614", "125×4=500")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
